$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second "Estado de Cuenta" data row (period 2506 for the same worker) -
# the whole row is deleted and everything below it shifts up one row.
$ws.Rows("17:17").Delete()

# Update the totals that summarize the remaining (single) period row:
# Valor Mora total (was sum of two 56940 rows = 113880, now just one row = 56940)
$ws.Range("E11").Value = 56940
# Cant. Periodos (was 2 periods, now only 1 remains)
$ws.Range("F13").Value = 1
